$d = $word.ActiveDocument

# 1) "4 Replies" -> "5 Replies" (hyperlink display text stays pointing to the
#    same #comments anchor; only the visible run text changes).
$d.Content.Find.Execute("4 Replies", $true, $true, $false, $false, $false,
                         $true, 1, $false, "5 Replies", 2)

# 2) Drop the trailing "This entry was posted in ..." blog-chrome block
#    (post categories/tags, post navigation, "N thoughts on", pingback list,
#    "Leave a Reply", and the Akismet notice) that ran from there through
#    the end of the document.
$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "This entry was posted in") {
        $startPara = $p
        break
    }
}

if ($startPara -ne $null) {
    $delRange = $d.Range($startPara.Range.Start, $d.Content.End)
    $delRange.Delete()
}

# 3) The paragraph left behind at the end of the document (previously just a
#    single blank-space run) gains two more blank-space runs.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$contentEnd = $r.End - 1
$targetRange = $d.Range($r.Start, $contentEnd)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$targetRange.InsertXML($xml)
